$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 19
$ws.Range("H19").Value = 1612
$ws.Range("I19").Value = 1583
$ws.Range("K19").Value = 1583
$ws.Range("M19").Value = -1408

# ALC row 80
$ws.Range("H80").Value = 2999
$ws.Range("J80").Value = 2999
$ws.Range("L80").Value = 8997
$ws.Range("N80").Value = -10993

# ALC row 83
$ws.Range("H83").Value = 2999
$ws.Range("J83").Value = 2999
$ws.Range("L83").Value = 26991
$ws.Range("N83").Value = -36975

# ALC row 98
$ws.Range("H98").Value = 2029.8889
$ws.Range("I98").Value = 1210.5
$ws.Range("K98").Value = 1210.5
$ws.Range("M98").Value = 287.5

# ALC row 100
$ws.Range("H100").Value = 3156.4285
$ws.Range("I100").Value = 3100
$ws.Range("J100").Value = 3231.6667
$ws.Range("K100").Value = 3100
$ws.Range("L100").Value = 3231.6667
$ws.Range("M100").Value = -2559
$ws.Range("N100").Value = -4313.6667

# ALC row 122
$ws.Range("H122").Value = 2029.8889
$ws.Range("I122").Value = 1210.5
$ws.Range("K122").Value = 3631.5
$ws.Range("M122").Value = -1181.5

$ws = $wb.Worksheets.Item("ARM")
# ARM row 96
$ws.Range("H96").Value = 11172
$ws.Range("J96").Value = 11172
$ws.Range("L96").Value = 11172
$ws.Range("N96").Value = -16664

# ARM row 110
$ws.Range("H110").Value = 687.9091
$ws.Range("I110").Value = 613.2857
$ws.Range("J110").Value = 818.5
$ws.Range("K110").Value = 613.2857
$ws.Range("L110").Value = 818.5
$ws.Range("M110").Value = 1431.7143
$ws.Range("N110").Value = -4908.5

# ARM row 124
$ws.Range("H124").Value = 19995
$ws.Range("J124").Value = 19995
$ws.Range("L124").Value = 19995
$ws.Range("N124").Value = -29815

$ws = $wb.Worksheets.Item("CRP")
# CRP row 10
$ws.Range("H10").Value = 151.33333
$ws.Range("I10").Value = 161.6
$ws.Range("J10").Value = 100
$ws.Range("K10").Value = 161.6
$ws.Range("L10").Value = 100
$ws.Range("M10").Value = -22.59999999999999
$ws.Range("N10").Value = -378

# CRP row 13
$ws.Range("H13").Value = 2000
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 2000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 2000
$ws.Range("N13").Value = -2278
$ws.Range("M13").ClearContents()

# CRP row 14
$ws.Range("H14").Value = 310
$ws.Range("J14").Value = 310
$ws.Range("L14").Value = 310
$ws.Range("M14").Value = -650

# CRP row 16
$ws.Range("H16").Value = 2005
$ws.Range("I16").Value = 2005
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2005
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1718
$ws.Range("N16").ClearContents()

# CRP row 21
$ws.Range("H21").Value = 26000
$ws.Range("J21").Value = 26000
$ws.Range("L21").Value = 26000
$ws.Range("N21").Value = -26470

# CRP row 68
$ws.Range("H68").Value = 39998.637
$ws.Range("J68").Value = 39998.637
$ws.Range("L68").Value = 39998.637
$ws.Range("N68").Value = -41496.637

# CRP row 71
$ws.Range("H71").Value = 39998.637
$ws.Range("J71").Value = 39998.637
$ws.Range("L71").Value = 119995.911
$ws.Range("N71").Value = -127483.911

# CRP row 92
$ws.Range("H92").Value = 45800
$ws.Range("J92").Value = 48250
$ws.Range("L92").Value = 48250
$ws.Range("N92").Value = -53242

# CRP row 113
$ws.Range("H113").Value = 2005
$ws.Range("I113").Value = 2005
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2005
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 165
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# CUL row 48
$ws.Range("H48").Value = 250
$ws.Range("I48").Value = 250
$ws.Range("K48").Value = 750
$ws.Range("M48").Value = -500

# CUL row 121
$ws.Range("H121").Value = 832.2857
$ws.Range("I121").Value = 641.1429000000001
$ws.Range("J121").Value = 1023.4286
$ws.Range("K121").Value = 1923.4287
$ws.Range("L121").Value = 3070.2858
$ws.Range("M121").Value = -613.4287000000002
$ws.Range("N121").Value = -5690.2858

# CUL row 131
$ws.Range("H131").Value = 810
$ws.Range("I131").Value = 810
$ws.Range("K131").Value = 2430
$ws.Range("M131").Value = 2610

$ws = $wb.Worksheets.Item("GSM")
# GSM row 43
$ws.Range("H43").Value = 15199
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 15199
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 15199
$ws.Range("N43").Value = -15501
$ws.Range("M43").ClearContents()

# GSM row 57
$ws.Range("H57").Value = 22395.8
$ws.Range("J57").Value = 24994.75
$ws.Range("L57").Value = 24994.75
$ws.Range("N57").Value = -26634.75

# GSM row 70
$ws.Range("H70").Value = 55562470
$ws.Range("I70").Value = 66672964
$ws.Range("K70").Value = 66672964
$ws.Range("M70").Value = -66672694

# GSM row 73
$ws.Range("H73").Value = 55562470
$ws.Range("I73").Value = 66672964
$ws.Range("K73").Value = 66672964
$ws.Range("M73").Value = -66672028

# GSM row 95
$ws.Range("H95").Value = 30000
$ws.Range("J95").Value = 30000
$ws.Range("L95").Value = 30000
$ws.Range("N95").Value = -35492

# GSM row 113
$ws.Range("H113").Value = 2000
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

# GSM row 132
$ws.Range("H132").Value = 3745.1765
$ws.Range("I132").Value = 3745.1765
$ws.Range("K132").Value = 11235.5295
$ws.Range("M132").Value = -8705.529500000001

$ws = $wb.Worksheets.Item("LTW")
# LTW row 68
$ws.Range("H68").Value = 52747.25
$ws.Range("J68").Value = 68999.336
$ws.Range("L68").Value = 68999.336
$ws.Range("N68").Value = -70497.336

# LTW row 71
$ws.Range("H71").Value = 52747.25
$ws.Range("J71").Value = 68999.336
$ws.Range("L71").Value = 344996.68
$ws.Range("N71").Value = -352484.68

# LTW row 132
$ws.Range("H132").Value = 2905.8462
$ws.Range("I132").Value = 2578.1
$ws.Range("K132").Value = 7734.299999999999
$ws.Range("M132").Value = -5204.299999999999

# LTW row 136
$ws.Range("H136").Value = 3444.8572
$ws.Range("I136").Value = 3268.1667
$ws.Range("J136").Value = 4505
$ws.Range("K136").Value = 9804.500100000001
$ws.Range("L136").Value = 13515
$ws.Range("M136").Value = -7254.500100000001
$ws.Range("N136").Value = -18615

$ws = $wb.Worksheets.Item("WVR")
# WVR row 122
$ws.Range("H122").Value = 2034.3182
$ws.Range("I122").Value = 1770.8889
$ws.Range("K122").Value = 5312.6667
$ws.Range("M122").Value = -2862.6667
